# Scheduled runner update: refresh market-price derived columns (H:N) on each leve-profit sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 684.625
$ws.Range("I2").Value = 684.625
$ws.Range("K2").Value = 684.625
$ws.Range("M2").Value = -571.625
$ws.Range("H17").Value = 3145715.5
$ws.Range("J17").Value = 3145715.5
$ws.Range("L17").Value = 9437146.5
$ws.Range("N17").Value = -9437482.5
$ws.Range("H31").Value = 131.66667
$ws.Range("I31").Value = 131.66667
$ws.Range("K31").Value = 395.00001
$ws.Range("M31").Value = -165.00001
$ws.Range("H74").Value = 5918.5
$ws.Range("I74").Value = 4834.3335
$ws.Range("K74").Value = 4834.3335
$ws.Range("M74").Value = -3898.3335
$ws.Range("H76").Value = 4046.25
$ws.Range("I76").Value = 3380.1333
$ws.Range("J76").Value = 5156.4443
$ws.Range("K76").Value = 3380.1333
$ws.Range("L76").Value = 5156.4443
$ws.Range("M76").Value = -3065.1333
$ws.Range("N76").Value = -5786.4443
$ws.Range("H77").Value = 5918.5
$ws.Range("I77").Value = 4834.3335
$ws.Range("K77").Value = 24171.6675
$ws.Range("M77").Value = -19491.6675
$ws.Range("H79").Value = 4046.25
$ws.Range("I79").Value = 3380.1333
$ws.Range("J79").Value = 5156.4443
$ws.Range("K79").Value = 3380.1333
$ws.Range("L79").Value = 5156.4443
$ws.Range("M79").Value = -2288.1333
$ws.Range("N79").Value = -7340.4443
$ws.Range("H138").Value = 14379.8
$ws.Range("I138").Value = 723.48486
$ws.Range("J138").Value = 22286.088
$ws.Range("K138").Value = 2170.45458
$ws.Range("L138").Value = 66858.264
$ws.Range("M138").Value = 2969.54542
$ws.Range("N138").Value = -77138.264

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 180
$ws.Range("I5").Value = 160
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 160
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -48
$ws.Range("N5").Value = -424
$ws.Range("H102").Value = 3089068.8
$ws.Range("J102").Value = 4666.6665
$ws.Range("L102").Value = 4666.6665
$ws.Range("N102").Value = -7910.6665
$ws.Range("H132").Value = 5988
$ws.Range("I132").Value = 1450.05
$ws.Range("J132").Value = 12969.462
$ws.Range("K132").Value = 4350.15
$ws.Range("L132").Value = 38908.386
$ws.Range("M132").Value = -1820.15
$ws.Range("N132").Value = -43968.386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 180
$ws.Range("I4").Value = 160
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 160
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -45
$ws.Range("N4").Value = -430
$ws.Range("H7").Value = 1453687.2
$ws.Range("I7").Value = 2507648.2
$ws.Range("J7").Value = 48406
$ws.Range("K7").Value = 2507648.2
$ws.Range("L7").Value = 48406
$ws.Range("M7").Value = -2507535.2
$ws.Range("N7").Value = -48632
$ws.Range("H22").Value = 165.73846
$ws.Range("I22").Value = 165.73846
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 165.73846
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 7.261539999999997
$ws.Range("N22").ClearContents()
$ws.Range("H94").Value = 1451.2325
$ws.Range("I94").Value = 1294
$ws.Range("J94").Value = 1908.6364
$ws.Range("K94").Value = 1294
$ws.Range("L94").Value = 1908.6364
$ws.Range("M94").Value = -843
$ws.Range("N94").Value = -2810.6364
$ws.Range("H99").Value = 111113336
$ws.Range("I99").Value = 500001000
$ws.Range("J99").Value = 2573.7144
$ws.Range("K99").Value = 500001000
$ws.Range("L99").Value = 2573.7144
$ws.Range("M99").Value = -499999502
$ws.Range("N99").Value = -5569.7144
$ws.Range("H105").Value = 2344.353
$ws.Range("I105").Value = 1979.9
$ws.Range("J105").Value = 2865
$ws.Range("K105").Value = 1979.9
$ws.Range("L105").Value = 2865
$ws.Range("M105").Value = -232.9000000000001
$ws.Range("N105").Value = -6359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7260.8687
$ws.Range("I31").Value = 1244.2273
$ws.Range("J31").Value = 15533.75
$ws.Range("K31").Value = 1244.2273
$ws.Range("L31").Value = 15533.75
$ws.Range("M31").Value = -949.2273
$ws.Range("N31").Value = -16123.75
$ws.Range("H34").Value = 7260.8687
$ws.Range("I34").Value = 1244.2273
$ws.Range("J34").Value = 15533.75
$ws.Range("K34").Value = 1244.2273
$ws.Range("L34").Value = 15533.75
$ws.Range("M34").Value = -1042.2273
$ws.Range("N34").Value = -15937.75
$ws.Range("H122").Value = 872
$ws.Range("I122").Value = 813.2857
$ws.Range("J122").Value = 1009
$ws.Range("K122").Value = 2439.8571
$ws.Range("L122").Value = 3027
$ws.Range("M122").Value = 10.14289999999983
$ws.Range("N122").Value = -7927
$ws.Range("H140").Value = 26087.375
$ws.Range("J140").Value = 26087.375
$ws.Range("L140").Value = 26087.375
$ws.Range("N140").Value = -36447.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 7142947
$ws.Range("I23").Value = 16666711
$ws.Range("J23").Value = 123.75
$ws.Range("K23").Value = 50000133
$ws.Range("L23").Value = 371.25
$ws.Range("M23").Value = -49999898
$ws.Range("N23").Value = -841.25
$ws.Range("H34").Value = 1808.2727
$ws.Range("I34").Value = 472.75
$ws.Range("J34").Value = 2571.4285
$ws.Range("K34").Value = 1418.25
$ws.Range("L34").Value = 7714.2855
$ws.Range("M34").Value = -1334.25
$ws.Range("N34").Value = -7882.2855
$ws.Range("H55").Value = 4372.143
$ws.Range("J55").Value = 4372.143
$ws.Range("L55").Value = 13116.429
$ws.Range("N55").Value = -13470.429
$ws.Range("H117").Value = 23811734
$ws.Range("I117").Value = 396.75
$ws.Range("J117").Value = 33336268
$ws.Range("K117").Value = 1190.25
$ws.Range("L117").Value = 100008804
$ws.Range("M117").Value = 2251.75
$ws.Range("N117").Value = -100015688
$ws.Range("H132").Value = 2199.9219
$ws.Range("I132").Value = 1920.25
$ws.Range("K132").Value = 17282.25
$ws.Range("M132").Value = -14752.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 54753
$ws.Range("J4").Value = 54753
$ws.Range("L4").Value = 54753
$ws.Range("N4").Value = -54977
$ws.Range("H70").Value = 5527.5557
$ws.Range("I70").Value = 5477.794
$ws.Range("J70").Value = 5681.364
$ws.Range("K70").Value = 5477.794
$ws.Range("L70").Value = 5681.364
$ws.Range("M70").Value = -5207.794
$ws.Range("N70").Value = -6221.364
$ws.Range("H73").Value = 5527.5557
$ws.Range("I73").Value = 5477.794
$ws.Range("J73").Value = 5681.364
$ws.Range("K73").Value = 5477.794
$ws.Range("L73").Value = 5681.364
$ws.Range("M73").Value = -4541.794
$ws.Range("N73").Value = -7553.364
$ws.Range("H80").Value = 2881.1738
$ws.Range("I80").Value = 2869.3635
$ws.Range("J80").Value = 2892
$ws.Range("K80").Value = 2869.3635
$ws.Range("L80").Value = 2892
$ws.Range("M80").Value = -1871.3635
$ws.Range("N80").Value = -4888
$ws.Range("H83").Value = 2881.1738
$ws.Range("I83").Value = 2869.3635
$ws.Range("J83").Value = 2892
$ws.Range("K83").Value = 14346.8175
$ws.Range("L83").Value = 14460
$ws.Range("M83").Value = -9354.817499999999
$ws.Range("N83").Value = -24444

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 16667808
$ws.Range("J46").Value = 1388.6666
$ws.Range("L46").Value = 1388.6666
$ws.Range("N46").Value = -1764.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 35936.93
$ws.Range("J46").Value = 35936.93
$ws.Range("L46").Value = 35936.93
$ws.Range("N46").Value = -36398.93
$ws.Range("H96").Value = 2395.75
$ws.Range("I96").Value = 1933
$ws.Range("K96").Value = 1933
$ws.Range("M96").Value = -560
$ws.Range("H107").Value = 76923440
$ws.Range("I107").Value = 111111420
$ws.Range("J107").Value = 475
$ws.Range("K107").Value = 333334260
$ws.Range("L107").Value = 1425
$ws.Range("M107").Value = -333332340
$ws.Range("N107").Value = -5265
$ws.Range("H126").Value = 1084.2727
$ws.Range("I126").Value = 843
$ws.Range("K126").Value = 2529
$ws.Range("M126").Value = -59
$ws.Range("H134").Value = 35936.93
$ws.Range("J134").Value = 35936.93
$ws.Range("L134").Value = 107810.79
$ws.Range("N134").Value = -112880.79
$ws.Range("H141").Value = 56123.89
$ws.Range("J141").Value = 60664.375
$ws.Range("L141").Value = 60664.375
$ws.Range("N141").Value = -71024.375
